$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$meta.Range("B9").Value = "Alvearie Team"

# Row 10: Contact/No display for ContactDetail -> Jurisdiction/United States of America
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" row - remove it entirely, shifting rows 12-21 up
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" (sheet2) ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element - Short/Definition updated to be specific to this profile
$elements.Range("K2").Value = "Stated Reason"
$elements.Range("L2").Value = "Recorded reason specified by the recipient"
